# Update the cryptos table (Price + Volume(1h) columns) with refreshed values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text is digits-and-a-single-dot would otherwise be auto-coerced
# to a number by plain Range.Value assignment (losing the original text formatting,
# e.g. trailing zeros). Route those through a temporary ="..." formula that is then
# flattened to a literal value via Copy/PasteSpecial(values-only), which keeps the
# cell a plain text value (same as the rest) without touching any cell formatting.
function Set-TextValue($cell, $text) {
    $rng = $ws.Range($cell)
    $rng.Formula = ('="' + $text + '"')
    $rng.Copy()
    $rng.PasteSpecial(-4163)
}
$excel.CutCopyMode = $false

$ws.Range("D2").Value = '67.225.10'
$ws.Range("E2").Value = '  +4.50%  '
$ws.Range("D3").Value = '3.260.17'
$ws.Range("E3").Value = '  +2.63%  '
$ws.Range("E4").Value = '  -0.04%  '
Set-TextValue 'D5' '578.19'
$ws.Range("E5").Value = '  +2.17%  '
Set-TextValue 'D6' '178.93'
$ws.Range("E6").Value = '  +4.98%  '
$ws.Range("E7").Value = '  -0.07%  '
Set-TextValue 'D8' '0.604'
$ws.Range("E8").Value = '  -0.94%  '
$ws.Range("D9").Value = '3.259.99'
$ws.Range("E9").Value = '  +2.69%  '
$ws.Range("E10").Value = '  +4.02%  '
$ws.Range("E11").Value = '  +1.75%  '
$ws.Range("E12").Value = '  +4.52%  '
$ws.Range("D13").Value = '3.826.21'
$ws.Range("E13").Value = '  +2.65%  '
$ws.Range("E14").Value = '  +0.74%  '
Set-TextValue 'D15' '28.27'
$ws.Range("E15").Value = '  +3.08%  '
$ws.Range("D16").Value = '67.186.71'
$ws.Range("E16").Value = '  +4.44%  '
$ws.Range("E17").Value = '  +2.78%  '
$ws.Range("D18").Value = '3.258.97'
$ws.Range("E18").Value = '  +2.62%  '
Set-TextValue 'D19' '5.88'
$ws.Range("E19").Value = '  +2.21%  '
Set-TextValue 'D20' '13.44'
$ws.Range("E20").Value = '  +3.47%  '
Set-TextValue 'D21' '374.78'
$ws.Range("E21").Value = '  +5.99%  '
$ws.Range("E22").Value = '  +6.40%  '
$ws.Range("E23").Value = '  -0.02%  '
Set-TextValue 'D24' '71.22'
$ws.Range("E24").Value = '  +3.12%  '
$ws.Range("E25").Value = '  +1.92%  '
$ws.Range("D26").Value = '3.396.61'
$ws.Range("E26").Value = '  +2.63%  '
$ws.Range("E27").Value = '  -0.78%  '
Set-TextValue 'D28' '9.88'
$ws.Range("E29").Value = '  +1.62%  '
$ws.Range("E30").Value = '  -0.31%  '
$ws.Range("E31").Value = '  +3.88%  '
$ws.Range("E32").Value = '  +0.42%  '
Set-TextValue 'D33' '22.63'
$ws.Range("E33").Value = '  +2.51%  '
$ws.Range("E34").Value = '  +0.08%  '
$ws.Range("E35").Value = '  +5.22%  '
$ws.Range("E36").Value = '  +2.61%  '
Set-TextValue 'D37' '167.26'
$ws.Range("E37").Value = '  +7.81%  '
$ws.Range("E38").Value = '  +4.57%  '
Set-TextValue 'D39' '0.858'
$ws.Range("E39").Value = '  +5.00%  '
$ws.Range("E40").Value = '  +10.32%  '
Set-TextValue 'D41' '27.07'
$ws.Range("E41").Value = '  +4.39%  '
Set-TextValue 'D42' '2.60'
$ws.Range("E42").Value = '  +1.49%  '
$ws.Range("D43").Value = '2.767.54'
$ws.Range("E43").Value = '  +5.55%  '
Set-TextValue 'D44' '6.51'
$ws.Range("E44").Value = '  +7.82%  '
Set-TextValue 'D45' '356.01'
$ws.Range("E45").Value = '  +10.07%  '
Set-TextValue 'D46' '4.41'
$ws.Range("E46").Value = '  +5.36%  '
Set-TextValue 'D47' '25.64'
$ws.Range("E47").Value = '  +7.35%  '
Set-TextValue 'D48' '40.50'
$ws.Range("E48").Value = '  +2.23%  '
Set-TextValue 'D49' '0.0676'
$ws.Range("E49").Value = '  +2.71%  '
Set-TextValue 'D50' '0.0280'
$ws.Range("E50").Value = '  +3.54%  '
$ws.Range("E51").Value = '  +0.70%  '
